$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.122.50'
$ws.Range('E2').Value = '  +3.85%  '
$ws.Range('D3').Value = '2.423.97'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('D5').Value = '553.38'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').Value = '138.79'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('E9').Value = '  +3.56%  '
$ws.Range('D10').Value = '5.79'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').Value = '25.00'
$ws.Range('E13').Value = '  +5.27%  '
$ws.Range('D14').Value = '2.855.24'
$ws.Range('E14').Value = '  +3.10%  '
$ws.Range('D15').Value = '60.048.95'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').Value = '0.0000139'
$ws.Range('E16').Value = '  +3.49%  '
$ws.Range('D17').Value = '2.421.64'
$ws.Range('E17').Value = '  +4.27%  '
$ws.Range('D18').Value = '11.38'
$ws.Range('E18').Value = '  +6.33%  '
$ws.Range('D19').Value = '4.40'
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('D20').Value = '332.84'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '6.79'
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = '65.21'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').Value = '8.62'
$ws.Range('E25').Value = '  +3.03%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '0.0₃0788'
$ws.Range('E28').Value = '  +7.24%  '
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('D30').Value = '6.31'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('D31').Value = '169.85'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +3.49%  '
$ws.Range('D33').Value = '18.72'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('E35').Value = '  +6.05%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '4.21'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.422'
$ws.Range('E39').Value = '  +11.53%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '321.84'
$ws.Range('E40').Value = '  +11.55%  '
$ws.Range('D41').Value = '39.49'
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('D42').Value = '3.70'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '139.71'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').Value = '0.0522'
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('D46').Value = '19.56'
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('D47').Value = '0.414'
$ws.Range('E47').Value = '  +9.33%  '
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('D50').Value = '17.81'
$ws.Range('E50').Value = '  +2.08%  '
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  -0.48%  '
